$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C width change: 12.7109375 -> 11.7109375 (XML character width units)
$ws.Columns.Item(3).ColumnWidth = 10.877604166666666

# Cell value updates
$ws.Range("F1").Value = 0.85498687024738396
$ws.Range("AF1").Value = 0.95393137257071392
$ws.Range("A2").Value = 0.78416164884345774
$ws.Range("Q2").Value = 0.83060336028936765
$ws.Range("E3").Value = 0.95076031997973831
$ws.Range("E4").Value = 0.7752390743550186
$ws.Range("AN4").Value = 0.78701823051812148
$ws.Range("Y5").Value = 0.77852446872220171
$ws.Range("W6").Value = 0.98441897553622604
$ws.Range("AU6").Value = 0.98436135691221893
$ws.Range("A8").Value = 0.90455708421368763
$ws.Range("E8").Value = 0.94395028574901585
$ws.Range("J8").Value = 0.9630396166848425
$ws.Range("AG8").Value = 0.96412352067550056
$ws.Range("AI9").Value = 0.74004659855184829
$ws.Range("BL9").Value = 0.56928162314664288
$ws.Range("L10").Value = 0.91644626966714338
$ws.Range("AJ10").Value = 0.81041112848726249
$ws.Range("V11").Value = 0.80550141107019146
$ws.Range("BK11").Value = 0.97434625753990756
$ws.Range("BP11").Value = 0.80916777540309459
$ws.Range("AO12").Value = 0.55218274519606725
$ws.Range("BE12").Value = 0.85153589361618387
$ws.Range("O13").Value = 0.70822846079035984
$ws.Range("BM13").Value = 0.84783964924025201
$ws.Range("S14").Value = 0.92353311475247946
$ws.Range("AD14").Value = 0.71326907200010425
$ws.Range("BB14").Value = 0.97940839268555013
$ws.Range("Z16").Value = 0.94052855168957772
$ws.Range("AK16").Value = 0.92244831255887427
$ws.Range("O17").Value = 0.60015879760849011
$ws.Range("BF17").Value = 0.66850452686531514
$ws.Range("O18").Value = 0.51345571401992451
$ws.Range("P18").Value = 0.82246155104953766
$ws.Range("Q18").Value = 0.90156783621490333
$ws.Range("W18").Value = 0.7906598303380139
$ws.Range("BD18").Value = 0.87863292442258079
$ws.Range("U19").Value = 0.63106337823734093
$ws.Range("Q20").Value = 0.60817377241995918
$ws.Range("BO20").Value = 0.89111444604301626
$ws.Range("B21").Value = 0.71375329044786828
$ws.Range("AC21").Value = 0.94992004783280826
$ws.Range("AD21").Value = 0.9483495904579744
$ws.Range("AY21").Value = 0.74225583815685958
$ws.Range("X22").Value = 0.64050075593266209
$ws.Range("BO22").Value = 0.96373873758364859
$ws.Range("AU23").Value = 0.6648004367351219
$ws.Range("BG23").Value = 0.59322550742158509
$ws.Range("AD24").Value = 0.7682732439946186
$ws.Range("BE24").Value = 0.93681558361451822
$ws.Range("AW25").Value = 0.99268268652857761
$ws.Range("AH26").Value = 0.64109067686866372
$ws.Range("AB27").Value = 0.94799239654929868
$ws.Range("BM27").Value = 0.83782269885619276
$ws.Range("D28").Value = 0.82909337285643259
$ws.Range("AR28").Value = 0.74288049906569664
$ws.Range("BB28").Value = 0.57615970351390922
$ws.Range("AQ29").Value = 0.98621053939804493
$ws.Range("BL29").Value = 0.79302223812095263
$ws.Range("D30").Value = 0.73156942890665999
$ws.Range("AN30").Value = 0.79564351452186521
$ws.Range("BM30").Value = 0.68842243349266985
$ws.Range("AG31").Value = 0.88476704080046087
$ws.Range("AT33").Value = 0.94688881912146661
$ws.Range("G34").Value = 0.92038371013466147
$ws.Range("AE34").Value = 0.76789243874520186
$ws.Range("BI34").Value = 0.74364389979701484
$ws.Range("S35").Value = 0.87178358316016369
$ws.Range("AT35").Value = 0.74190367593339823
$ws.Range("AC37").Value = 0.71767562052631195
$ws.Range("I38").Value = 0.62067690161346567
$ws.Range("BP38").Value = 0.91391141285320521
$ws.Range("T39").Value = 0.74947547423915029
$ws.Range("BI39").Value = 0.8425806142578347
$ws.Range("AM40").Value = 0.95814603719574243
$ws.Range("BJ40").Value = 0.68621793457018776
$ws.Range("AM41").Value = 0.85721532823989377
$ws.Range("AQ41").Value = 0.90893404557903157
$ws.Range("AV42").Value = 0.8635278694658024
$ws.Range("BO42").Value = 0.79008653370281023
$ws.Range("I44").Value = 0.67874347656035394
$ws.Range("J45").Value = 0.52824894378562259
$ws.Range("Q45").Value = 0.6633203061814168
$ws.Range("AT45").Value = 0.81335052316286394
$ws.Range("L46").Value = 0.91080694773838555
$ws.Range("C47").Value = 0.99919766197086968
$ws.Range("AE47").Value = 0.99168118790609783
$ws.Range("BH47").Value = 0.6602683391715426
$ws.Range("BL48").Value = 0.90518136487794498
$ws.Range("AC49").Value = 0.75017254853523463
$ws.Range("AE49").Value = 0.74338809681612128
$ws.Range("AW50").Value = 0.96108707625798906
$ws.Range("AZ50").Value = 0.8595131695709477
$ws.Range("BA50").Value = 0.85253074842656362
$ws.Range("BB50").Value = 0.86987683554560269
$ws.Range("BE51").Value = 0.97993151221041896
$ws.Range("BI51").Value = 0.66884451105271148
$ws.Range("BK51").Value = 0.93820250355220702
$ws.Range("P52").Value = 0.76134369112401323
$ws.Range("AO52").Value = 0.92973925012031533
$ws.Range("O53").Value = 0.83494007436054107
$ws.Range("AT53").Value = 0.70320467825405331
$ws.Range("BM53").Value = 0.87614785646320215
$ws.Range("U54").Value = 0.97445474623880235
$ws.Range("M55").Value = 0.92138963184888878
$ws.Range("W55").Value = 0.70289120352199452
$ws.Range("AF55").Value = 0.80187127426857052
$ws.Range("AG55").Value = 0.63215854563025697
$ws.Range("AR55").Value = 0.71257846637760747
$ws.Range("AT55").Value = 0.8660792147924451
$ws.Range("BI55").Value = 0.77613683286963164
$ws.Range("AV56").Value = 0.99925294392696873
$ws.Range("BM56").Value = 0.88792154158262737
$ws.Range("W57").Value = 0.90718493640478304
$ws.Range("P58").Value = 0.83653713290032194
$ws.Range("BD58").Value = 0.89629641488710643
$ws.Range("BF59").Value = 0.60368194740157533
$ws.Range("BI59").Value = 0.91395192603652464
$ws.Range("G60").Value = 0.92644377959561652
$ws.Range("R60").Value = 0.6461475700869741
$ws.Range("AI60").Value = 0.93869573945314244
$ws.Range("AQ60").Value = 0.99291131254780174
$ws.Range("BI60").Value = 0.96875790469893075
$ws.Range("F61").Value = 0.61115377231529799
$ws.Range("AS61").Value = 0.91942276243995069
$ws.Range("BL62").Value = 0.76694240168852035
$ws.Range("AJ63").Value = 0.81118862667380376
$ws.Range("AP63").Value = 0.70646068707827092
$ws.Range("E64").Value = 0.92486980175470046
$ws.Range("BI64").Value = 0.93611472931350292
$ws.Range("V65").Value = 0.89607570401388303
$ws.Range("BF65").Value = 0.97741174280315568
$ws.Range("AA66").Value = 0.72485347831311686
$ws.Range("BI66").Value = 0.90873892355127617
$ws.Range("AB67").Value = 0.93172877985816394
$ws.Range("BC68").Value = 0.86807420133542301
